$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("expert-selected")
$ws.Range("A2").Value = "no"
